$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.7281672816728167
$ws.Range("C2").Value = 0.4575278265963679
$ws.Range("D2").Value = 0.8135416666666667
$ws.Range("E2").Value = 0.5856767904011998
$ws.Range("F2").Value = 0.4434542885967195
$ws.Range("G2").Value = 0.4061507830270435
$ws.Range("H2").Value = 0.7576565016103061
